# MRD-1840: Update Part A template to retain form field.
# Converts the plain-text "{{ppcs_query_emails}}" / "{{revocation_order_recipients}}"
# merge placeholders into Word FORMTEXT form fields, so the legacy field survives
# until the probationAdmin flag supplies the real data.

$d = $word.ActiveDocument
$wdFindContinue = 1

function Clear-RangeText($rng) {
    $rng.Text = ""
}

# ---------------------------------------------------------------------------
# Hunk 1: first "{{ppcs_query_emails}}" (the standalone placeholder paragraph
# right under "E-mail address to which PPCS should respond ... mailbox:")
# ---------------------------------------------------------------------------
$scopeEnd = $d.Content.End
$firstScope = $d.Range(0, $scopeEnd)
$found = $firstScope.Find.Execute("{{ppcs_query_emails}}", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if (-not $found) {
    throw "Hunk1: could not find first {{ppcs_query_emails}} placeholder"
}
$target = $d.Range($firstScope.Start, $firstScope.End)
Clear-RangeText $target

$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:fldChar w:fldCharType="begin"><w:ffData><w:name w:val=""/><w:enabled/><w:calcOnExit w:val="0"/><w:statusText w:type="text" w:val="senior manager authorisation  - e mail address:"/><w:textInput><w:default w:val="{{countersign_aco_email}}"/></w:textInput></w:ffData></w:fldChar></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:instrText xml:space="preserve"> FORMTEXT </w:instrText></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:noProof/></w:rPr><w:t>{{ppcs_query_emails}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:noProof/></w:rPr><w:t xml:space="preserve">}      </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@
$target.InsertXML($xml1)

Write-Output "Hunk1 done"

# ---------------------------------------------------------------------------
# Hunk 2: second "{{ppcs_query_emails}}" -- this one shares a paragraph with
# the preceding "E-mail address ... mailbox:" label and a leading space, and
# must be split into its own paragraph (cloning the bordered pPr) before the
# placeholder text is turned into the FORMTEXT field.
# ---------------------------------------------------------------------------
$full = $d.Content
$found2 = $full.Find.Execute(": {{ppcs_query_emails}}", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if (-not $found2) {
    throw "Hunk2: could not find ': {{ppcs_query_emails}}' context"
}
$colonPos = $full.Start + 1
$splitPoint = $d.Range($colonPos, $colonPos)
$splitPoint.InsertParagraphAfter()

# Re-find the placeholder text; it now lives alone (with a leading space) in
# its own paragraph.
$scope2 = $d.Content
$found3 = $scope2.Find.Execute("{{ppcs_query_emails}}", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if (-not $found3) {
    throw "Hunk2: could not re-find {{ppcs_query_emails}} after paragraph split"
}
$newPara = $scope2.Paragraphs(1)
$paraRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
Clear-RangeText $paraRange

$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:fldChar w:fldCharType="begin"><w:ffData><w:name w:val=""/><w:enabled/><w:calcOnExit w:val="0"/><w:textInput/></w:ffData></w:fldChar></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:instrText xml:space="preserve"> FORMTEXT </w:instrText></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>{{ppcs_query_emails}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">}       </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@
$paraRange.InsertXML($xml2)

Write-Output "Hunk2 done"
